$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 2, shifting the existing data (rows 2-21) down to rows 3-22.
$ws.Rows.Item(2).Insert()

# Populate the newly inserted row 2 with the new accelerometer sample.
$ws.Range("A2").Value = 7.345117568969727
$ws.Range("B2").Value = -12.58289909362793
$ws.Range("C2").Value = 3.90805721282959

# Append the 9 new rows of accelerometer samples at the bottom (rows 23-31).
$newRows = @(
    @(-8.928971290588379, -17.86810111999512, 8.281005859375),
    @(-14.49608421325684, -1.527808666229248, 44.4189453125),
    @(-12.06443023681641, 6.844409942626953, 19.9449577331543),
    @(6.954762935638428, -76.15243530273438, 24.18494606018066),
    @(6.384909629821777, 5.00542688369751, -29.23712921142578),
    @(-34.79932403564453, -7.816071510314941, 1.089200496673584),
    @(-17.0820198059082, -31.8654670715332, 12.90904235839844),
    @(2.159783363342285, 0.4922776222229004, 7.778494358062744),
    @(3.85674524307251, 1.991205930709839, 21.4826774597168)
)

$row = 23
foreach ($values in $newRows) {
    $ws.Range("A$row").Value = $values[0]
    $ws.Range("B$row").Value = $values[1]
    $ws.Range("C$row").Value = $values[2]
    $row++
}
